$d = $word.ActiveDocument

# Update the date line
$d.Content.Find.Execute("2025-07-06 Sunday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-07-07 Monday", 2) | Out-Null

# Update the division-problem table. Data rows are 1, 5, 9, 13, 17 (the
# other rows are blank spacer rows), each with 5 columns. Addressed by
# position (row/col) rather than text search because several cells share
# identical source text ("36÷6=6, 0" appears twice) but map to different
# targets.
$t = $d.Tables.Item(1)

$values = @(
    @(1, 1, "26÷2=13, 0"),
    @(1, 2, "30÷6=5, 0"),
    @(1, 3, "16÷4=4, 0"),
    @(1, 4, "76÷7=10, 6"),
    @(1, 5, "42÷9=4, 6"),

    @(5, 1, "63÷6=10, 3"),
    @(5, 2, "83÷2=41, 1"),
    @(5, 3, "49÷9=5, 4"),
    @(5, 4, "41÷2=20, 1"),
    @(5, 5, "86÷8=10, 6"),

    @(9, 1, "84÷4=21, 0"),
    @(9, 2, "41÷5=8, 1"),
    @(9, 3, "77÷5=15, 2"),
    @(9, 4, "10÷4=2, 2"),
    @(9, 5, "86÷6=14, 2"),

    @(13, 1, "61÷8=7, 5"),
    @(13, 2, "43÷5=8, 3"),
    @(13, 3, "27÷9=3, 0"),
    @(13, 4, "69÷8=8, 5"),
    @(13, 5, "40÷4=10, 0"),

    @(17, 1, "55÷4=13, 3"),
    @(17, 2, "63÷8=7, 7"),
    @(17, 3, "20÷8=2, 4"),
    @(17, 4, "89÷6=14, 5"),
    @(17, 5, "42÷7=6, 0")
)

foreach ($entry in $values) {
    $row = $entry[0]
    $col = $entry[1]
    $text = $entry[2]
    $t.Cell($row, $col).Range.Text = $text
}
